$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- RUBRIC block (rows 7-10): fill in explicit zeros for previously-blank cells ---
# Row 8 - Exercises
$ws.Range("F8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("K8").Value = 0

# Row 9 - Projects
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0

# --- Clear the old ACTUAL / ACTUAL 2 blocks so we can rebuild them cleanly ---
$ws.Range("B11:K27").ClearContents()

# --- ACTUAL block (week 7 actuals), now starting at row 12 ---
$ws.Range("B12").Value = "ACTUAL"

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 6
$ws.Range("J13").Value = 7
$ws.Range("K13").Value = 8

$ws.Range("B14").Value = "Discussion"
$ws.Range("C14").Formula = "=SUM(D14:K14)"
$ws.Range("D14").Value = 1.5
$ws.Range("E14").Value = 0.75
$ws.Range("F14").Value = 0.75
$ws.Range("G14").Formula = "=F14"
$ws.Range("H14").Formula = "=G14"
$ws.Range("I14").Formula = "=H14"
$ws.Range("J14").Formula = "=I14"
$ws.Range("K14").Formula = "=J14"

$ws.Range("B15").Value = "Exercises"
$ws.Range("C15").Formula = "=SUM(D15:K15)"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 9.5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 0

$ws.Range("B16").Value = "Projects"
$ws.Range("C16").Formula = "=SUM(D16:K16)"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 14.25
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 14.7
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 7.5

$ws.Range("C17").Formula = "=SUM(C14:C16)"

# --- ACTUAL 2 block (week 8 actuals so far), now starting at row 22 ---
$ws.Range("B22").Value = "ACTUAL 2"

$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 7
$ws.Range("K23").Value = 8

$ws.Range("B24").Value = "Discussion"
$ws.Range("C24").Formula = "=SUM(D24:K24)"
$ws.Range("D24").Value = 1.5
$ws.Range("E24").Value = 0.75
$ws.Range("F24").Value = 0.75
$ws.Range("G24").Value = 0.75
$ws.Range("H24").Value = 0.75
$ws.Range("I24").Value = 1.5
$ws.Range("J24").Value = 1.5
$ws.Range("K24").Value = 1.5

$ws.Range("B25").Value = "Exercises"
$ws.Range("C25").Formula = "=SUM(D25:K25)"
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 10
$ws.Range("G25").Value = 9.5
$ws.Range("I25").Value = 9
$ws.Range("J25").Value = 9

$ws.Range("B26").Value = "Projects"
$ws.Range("C26").Formula = "=SUM(D26:K26)"
$ws.Range("F26").Value = 14.25
$ws.Range("H26").Value = 10
$ws.Range("K26").Value = 10

$ws.Range("C27").Formula = "=SUM(C24:C26)"

# --- Update the view: scroll back to top-left and select K17 ---
[void]$ws.Range("A1").Select()
[void]$ws.Range("K17").Select()
